# Updates to fit workbook.
# - Fill in "ready to be fit" markers for dataset row 9 (HighKick) and row 10 (Endgame)
#   in the newly added columns U and V.
# - Row 11 (9d dataset): clear the "ready to be fit" marker in column C (just keep the
#   green background), add a blank green-highlighted cell in column F, and mark
#   columns U ("hadding in 01") and V ("ready to be fit") as done.
# - Move the active selection from E14 to E16.
# - Add a new note under the "Other things..." section: "3500 MeV cut for pileup plots".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors matching the existing conditional fills used throughout the sheet
$greenFill  = 5296274   # RGB 92D050 (OLE/BGR encoded) -> used for style s=1
$orangeFill = 49407     # RGB FFC000 (OLE/BGR encoded) -> used for style s=2

# --- Row 9 (HighKick): add "ready to be fit" marks in U9 and V9 ---
$ws.Range("U9").Value = "ready to be fit"
$ws.Range("U9").Interior.Color = $orangeFill
$ws.Range("V9").Value = "ready to be fit"
$ws.Range("V9").Interior.Color = $orangeFill

# --- Row 10 (Endgame): add "ready to be fit" marks in U10 and V10 ---
$ws.Range("U10").Value = "ready to be fit"
$ws.Range("U10").Interior.Color = $orangeFill
$ws.Range("V10").Value = "ready to be fit"
$ws.Range("V10").Interior.Color = $orangeFill

# --- Row 11 (9d): clear C11's "ready to be fit" text but keep it highlighted green ---
$ws.Range("C11").ClearContents()
$ws.Range("C11").Interior.Color = $greenFill

# Add a blank green cell in F11
$ws.Range("F11").Interior.Color = $greenFill

# --- Add new to-do note under "Other things to do/ think about/ argue away" ---
# (added before "hadding in 01" so new shared strings line up with the source order)
$ws.Range("A23").Value = "3500 MeV cut for pileup plots"

# Mark U11 ("hadding in 01") and V11 ("ready to be fit") as done
$ws.Range("U11").Value = "hadding in 01"
$ws.Range("U11").Interior.Color = $orangeFill
$ws.Range("V11").Value = "ready to be fit"
$ws.Range("V11").Interior.Color = $orangeFill

# --- Move the active selection ---
$ws.Range("E16").Select()
